# Update mods data [2026-02-10 15:54:24]
# Appends a new daily row (row 92) to the ModCounts sheet:
#   A92 = "2026/02/10" (text, not a date), B92 = "逃离鸭科夫", C92 = 1183
# The new row reuses the same cell style (centered alignment) as the
# preceding data rows.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Seed row 92 with the same formatting (style) as the last existing data
# row (91) by copying its formats down first.
$ws.Range("A91:C91").Copy()
$ws.Range("A92:C92").PasteSpecial(-4122)

# Write the new values. The date-looking string in column A must stay text
# (matching how the log was originally written), so it is entered with a
# leading apostrophe to suppress Excel's automatic date conversion.
$ws.Cells.Item(92, 1).Value = "'2026/02/10"
$ws.Cells.Item(92, 2).Value = "逃离鸭科夫"
$ws.Cells.Item(92, 3).Value = 1183

# Re-apply the row-91 formatting once more so the forced-text cell in
# column A ends up sharing the exact same style as columns B and C (and
# as the rest of the table) instead of getting its own ad-hoc style.
$ws.Range("A91:C91").Copy()
$ws.Range("A92:C92").PasteSpecial(-4122)
